$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "A"
$ws.Range("B1").Value = "B"
$ws.Range("C1").Value = "C"
$ws.Range("D1").Value = "D"
$ws.Range("E1").Value = "E"
$ws.Range("F1").Value = "F"

# Row 2: time formatted with colons
$ws.Range("A2:F2").Value = "1:4:45"

# Row 3: time formatted with pipes
$ws.Range("A3:F3").Value = "1|4|45"

# Row 4: time formatted with spaces
$ws.Range("A4:F4").Value = "1 4 45"

# Row 5: time formatted with dashes
$ws.Range("A5:F5").NumberFormat = "@"
$ws.Range("A5:F5").Value = "1-4-45"
